$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New patient rows (3-5) being added. Pre-format the target range as Text so that
# numeric-looking values (ages, dates, lab results, etc.) are stored as literal
# strings, matching the original template row formatting, then restore the default
# (Normal) style so no stray number-format/style index is left on the cells.
$dataRange = $ws.Range("A3:AR5")
$dataRange.NumberFormat = "@"

# Row 3
$ws.Cells.Item(3, 1).Value = "Adriana Hernandez Caballero"
$ws.Cells.Item(3, 2).Value = "M"
$ws.Cells.Item(3, 3).Value = "1978-09-07"
$ws.Cells.Item(3, 4).Value = "40"
$ws.Cells.Item(3, 5).Value = "2019-04-01"
$ws.Cells.Item(3, 6).Value = "168"
$ws.Cells.Item(3, 7).Value = "58"
$ws.Cells.Item(3, 8).Value = "20"
$ws.Cells.Item(3, 9).Value = "1.65 m²"
$ws.Cells.Item(3, 10).Value = "8"
$ws.Cells.Item(3, 11).Value = "7"
$ws.Cells.Item(3, 12).Value = "2"
$ws.Cells.Item(3, 13).Value = "No"
$ws.Cells.Item(3, 14).Value = "0"
$ws.Cells.Item(3, 15).Value = "No"
$ws.Cells.Item(3, 16).Value = "5"
$ws.Cells.Item(3, 17).Value = "No"
$ws.Cells.Item(3, 18).Value = "10"
$ws.Cells.Item(3, 19).Value = "12"
$ws.Cells.Item(3, 20).Value = "5"
$ws.Cells.Item(3, 21).Value = "7"
$ws.Cells.Item(3, 22).Value = "15"
$ws.Cells.Item(3, 23).Value = "IgG"
$ws.Cells.Item(3, 24).Value = "Kappa"
$ws.Cells.Item(3, 25).Value = "2"
$ws.Cells.Item(3, 26).Value = "6"
$ws.Cells.Item(3, 27).Value = "4"
$ws.Cells.Item(3, 28).Value = "8"
$ws.Cells.Item(3, 29).Value = "No"
$ws.Cells.Item(3, 30).Value = "No"
$ws.Cells.Item(3, 31).Value = "5"
$ws.Cells.Item(3, 32).Value = "4"
$ws.Cells.Item(3, 33).Value = "8"
$ws.Cells.Item(3, 34).Value = "15"
$ws.Cells.Item(3, 35).Value = "9"
$ws.Cells.Item(3, 36).Value = "5"
$ws.Cells.Item(3, 37).Value = "25"
$ws.Cells.Item(3, 38).Value = "10"
$ws.Cells.Item(3, 39).Value = "22"
$ws.Cells.Item(3, 40).Value = "10"
$ws.Cells.Item(3, 41).Value = "No"
$ws.Cells.Item(3, 42).Value = "No"
$ws.Cells.Item(3, 43).Value = "No"

# Row 4
$ws.Cells.Item(4, 1).Value = "Julio Juarez Mendoza"
$ws.Cells.Item(4, 2).Value = "H"
$ws.Cells.Item(4, 3).Value = "1984-06-14"
$ws.Cells.Item(4, 4).Value = "35"
$ws.Cells.Item(4, 5).Value = "2017-06-08"
$ws.Cells.Item(4, 6).Value = "180"
$ws.Cells.Item(4, 7).Value = "88"
$ws.Cells.Item(4, 8).Value = "27"
$ws.Cells.Item(4, 9).Value = "2.1 m²"
$ws.Cells.Item(4, 10).Value = "9"
$ws.Cells.Item(4, 11).Value = "7"
$ws.Cells.Item(4, 12).Value = "1"
$ws.Cells.Item(4, 13).Value = "Si"
$ws.Cells.Item(4, 14).Value = "0"
$ws.Cells.Item(4, 15).Value = "Si"
$ws.Cells.Item(4, 16).Value = "5"
$ws.Cells.Item(4, 17).Value = "No"
$ws.Cells.Item(4, 18).Value = "10"
$ws.Cells.Item(4, 19).Value = "8"
$ws.Cells.Item(4, 20).Value = "10"
$ws.Cells.Item(4, 21).Value = "12"
$ws.Cells.Item(4, 22).Value = "15"
$ws.Cells.Item(4, 23).Value = "IgG"
$ws.Cells.Item(4, 24).Value = "No secretor"
$ws.Cells.Item(4, 25).Value = "16"
$ws.Cells.Item(4, 26).Value = "6"
$ws.Cells.Item(4, 27).Value = "4"
$ws.Cells.Item(4, 28).Value = "23"
$ws.Cells.Item(4, 29).Value = "No"
$ws.Cells.Item(4, 30).Value = "No"
$ws.Cells.Item(4, 31).Value = "8"
$ws.Cells.Item(4, 32).Value = "5"
$ws.Cells.Item(4, 33).Value = "9"
$ws.Cells.Item(4, 34).Value = "15"
$ws.Cells.Item(4, 35).Value = "9"
$ws.Cells.Item(4, 36).Value = "10"
$ws.Cells.Item(4, 37).Value = "25"
$ws.Cells.Item(4, 38).Value = "10"
$ws.Cells.Item(4, 39).Value = "22"
$ws.Cells.Item(4, 40).Value = "10"
$ws.Cells.Item(4, 41).Value = "No"
$ws.Cells.Item(4, 42).Value = "Si"
$ws.Cells.Item(4, 43).Value = "Si"

# Row 5
$ws.Cells.Item(5, 1).Value = "Eduardo Rosas Lopez"
$ws.Cells.Item(5, 2).Value = "H"
$ws.Cells.Item(5, 3).Value = "1985-06-04"
$ws.Cells.Item(5, 4).Value = "25"
$ws.Cells.Item(5, 5).Value = "2010-06-08"
$ws.Cells.Item(5, 6).Value = "180"
$ws.Cells.Item(5, 7).Value = "98"
$ws.Cells.Item(5, 8).Value = "30"
$ws.Cells.Item(5, 9).Value = "2.21 m²"
$ws.Cells.Item(5, 10).Value = "9"
$ws.Cells.Item(5, 11).Value = "7"
$ws.Cells.Item(5, 12).Value = "2"
$ws.Cells.Item(5, 13).Value = "No"
$ws.Cells.Item(5, 14).Value = "0"
$ws.Cells.Item(5, 15).Value = "Si"
$ws.Cells.Item(5, 16).Value = "5"
$ws.Cells.Item(5, 17).Value = "No"
$ws.Cells.Item(5, 18).Value = "15"
$ws.Cells.Item(5, 19).Value = "4"
$ws.Cells.Item(5, 20).Value = "5"
$ws.Cells.Item(5, 21).Value = "5"
$ws.Cells.Item(5, 22).Value = "16"
$ws.Cells.Item(5, 23).Value = "IgG"
$ws.Cells.Item(5, 24).Value = "No secretor"
$ws.Cells.Item(5, 25).Value = "16"
$ws.Cells.Item(5, 26).Value = "14"
$ws.Cells.Item(5, 27).Value = "4"
$ws.Cells.Item(5, 28).Value = "8"
$ws.Cells.Item(5, 29).Value = "No"
$ws.Cells.Item(5, 30).Value = "No"
$ws.Cells.Item(5, 31).Value = "8"
$ws.Cells.Item(5, 32).Value = "6"
$ws.Cells.Item(5, 33).Value = "9"
$ws.Cells.Item(5, 34).Value = "15"
$ws.Cells.Item(5, 35).Value = "9"
$ws.Cells.Item(5, 36).Value = "21"
$ws.Cells.Item(5, 37).Value = "25"
$ws.Cells.Item(5, 38).Value = "10"
$ws.Cells.Item(5, 39).Value = "22"
$ws.Cells.Item(5, 40).Value = "10"
$ws.Cells.Item(5, 41).Value = "No"
$ws.Cells.Item(5, 42).Value = "No"
$ws.Cells.Item(5, 43).Value = "Si"

# Restore default styling (removes the temporary text-number-format/quote-prefix
# styling so the new cells have no explicit style, same as the target workbook).
$dataRange.Style = "Normal"